# Session-1 follow-up edits to the Portland Housing Prices workbook.
$wb = $excel.ActiveWorkbook

# --- Workbook window position tweak ------------------------------------
$wb.Windows.Item(1).Left = 1940

# --- "Data" sheet updates ------------------------------------------------
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Relabel the existing headers and add a new "Prediction" column.
$ws.Range("A1").Value = "Num Bedrooms(x1)"
$ws.Range("B1").Value = "Area (Square Feet) (x2)"
$ws.Range("C1").Value = "Price ($) - output or target"
$ws.Range("D1").Value = "Prediction"

# Match the bold header styling used by the rest of row 1.
$ws.Range("D1").Font.Bold = $true

# Add a COST label below the data table.
$ws.Range("C49").Value = "COST"

# Zoom in and freeze the header row.
$win = $ws.Application.ActiveWindow
$win.Zoom = 116
$ws.Range("A2").Select()
$win.FreezePanes = $true
